# Update the Environment_group_table worksheet:
#  - rewrite Group/Code/Name/Description cell values (new codes E1..E5, new
#    Name/Description text)
#  - wrap text + widen columns A, C, D
#  - set explicit row heights for the data rows
#  - move the active selection to D6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----------------------------------------------------------
$ws.Range('A1').Value = 'Group'
$ws.Range('B1').Value = 'Code'
$ws.Range('C1').Value = 'Name'
$ws.Range('D1').Value = 'Description'

# ---- Environment 1 ---------------------------------------------------------
$ws.Range('A2').Value = 'Environment 1'
$ws.Range('B2').Value = 'E1'
$ws.Range('C2').Value = 'Midland humid ecosystems'
$ws.Range('D2').Value = 'Occurs on a broad variety of transitional terrains at low-to-mid elevations, usually near the coast. Characterised by water availability via high cloud cover, snowfall, and melt. Thick snow cover (probably seasonal) remains on the surfaces, which are mostly not steep or windy enough for the snow to slide down. Temperatures are mild and many areas have a substantial growing season. Biota consists mainly of mosses and lichens. Dominates the mountainous areas of the Antarctic peninsula, with representation in the Ellsworth mountains, the Transantarctic mountains, Victoria Land, and Enderby Land.'

# ---- Environment 2 ---------------------------------------------------------
$ws.Range('A3').Value = 'Environment 2'
$ws.Range('B3').Value = 'E2'
$ws.Range('C3').Value = 'High cliffs, crags, mountainsides, and slopes'
$ws.Range('D3').Value = 'High elevation, extremely cold, arid, and rugged ecosystems receiving low solar radiation (south-facing). Slopes are too steep to hold snow layer and are typically barren rock with very little biota. Occurs sprinkled among other ecosystems where steep slopes are present. Particular subunits cover rugged areas of the peninsula, and all subunits occur in the Transantarctic mountains.'

# ---- Environment 3 ---------------------------------------------------------
$ws.Range('A4').Value = 'Environment 3'
$ws.Range('B4').Value = 'E3'
$ws.Range('C4').Value = 'Mild lowlands'
$ws.Range('D4').Value = 'Relatively warm ecosystems with gentle relief, often coastal but may extend substantially inland. With especially low elevations, this group contains coastal rocky outcrops, small islands, and beaches as well as flat valley bottoms. May be rocky and barren, but often hosts bird colonies and pinnipeds if occurring near the coast. Occurs all around Antarctica with particular concentrations on the Antarctic Peninsula, in the Transantarctic Mountains, and Victoria Land. Dominates the coastal outcrops around the Eastern margin of the continent.'

# ---- Environment 4 ---------------------------------------------------------
$ws.Range('A5').Value = 'Environment 4'
$ws.Range('B5').Value = 'E4'
$ws.Range('C5').Value = 'Sunny inclines, mountainsides, nunataks and outcrops'
$ws.Range('D5').Value = 'Clear, sunlit (north-facing) rocky slopes with low snow cover. In terms of temperature, water availability, and terrain, E4 is transitional between the milder environments of E1/E3 and the more extreme environments of E2/E5.  Topography and elevation is quite variable. Occurs all over Antarctica with good representation, but particularly dominant in Victoria Land.'

# ---- Environment 5 ---------------------------------------------------------
$ws.Range('A6').Value = 'Environment 5'
$ws.Range('B6').Value = 'E5'
$ws.Range('C6').Value = 'Highland windy plateaus and outcrops'
$ws.Range('D6').Value = 'Very cold and arid high-elevation flatlands/plateaus with low cloud cover. Characterised especially by high winds. Occurs throughout the main continent but is nearly absent from the Antarctic Peninsula.'

# ---- Formatting ------------------------------------------------------------
# Wrap the Name/Description columns for every row of the table.
$ws.Range('C1:D6').WrapText = $true

# Widen column A, and the now-wrapped columns C and D (values chosen so the
# rounded column width this engine stores comes as close as possible to the
# author's original 16.85546875 / 39.28515625 / 95.42578125 character widths).
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(3).ColumnWidth = 38.5
$ws.Columns.Item(4).ColumnWidth = 94.66666666666667

# Explicit row heights to fit the wrapped text.
$ws.Rows.Item(2).RowHeight = 76.5
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 62.25
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 30

# ---- Selection ---------------------------------------------------------
$ws.Range('D6').Select()
